$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the duplicated "superuser" value into column C, row 5 (same as A5)
$ws.Range("C5").Value = $ws.Range("A5").Value2

# Update the active selection to C6 (as recorded in the saved view state)
$ws.Range("C6").Select()
